$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.478.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.59%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.096.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "582.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.088.55"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.84%  "

$ws.Range("E9").Value = "  -0.51%  "

$ws.Range("E10").Value = "  +5.96%  "

$ws.Range("E11").Value = "  -3.07%  "

$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.21%  "

$ws.Range("E15").Value = "  -1.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.605.25"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.88%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.322.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.094.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.57%  "

$ws.Range("E21").Value = "  +0.67%  "

$ws.Range("E23").Value = "  -1.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.90"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "  -2.17%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.31%  "

$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.67"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("E31").Value = "  -2.14%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.02%  "

$ws.Range("E33").Value = "  -1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.62"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0851"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.95%  "

$ws.Range("E36").Value = "  +2.56%  "

$ws.Range("E37").Value = "  -1.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.30"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.11%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.30%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "433.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0367"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.881.20"
$ws.Range("D44").Style = "Normal"

$ws.Range("E45").Value = "  -3.60%  "

$ws.Range("E46").Value = "  -3.58%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "35.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.67%  "

$ws.Range("E49").Value = "  -0.09%  "

$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.05"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.51%  "
